# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Membrillo" (Vega Modelo de Temuco)
# at the top of this product's data block (row 177), pushing the existing
# rows 177-212 down to 178-213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 177..212 down to 178..213 and open up a fresh row 177.
$ws.Rows.Item(177).Insert()

# Populate the new row 177 with this week's record.
$ws.Cells.Item(177, 1).Value  = 10
$ws.Cells.Item(177, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(177, 3).Value  = "La Araucanía"
$ws.Cells.Item(177, 4).Value  = 44782
$ws.Cells.Item(177, 5).Value  = 9
$ws.Cells.Item(177, 6).Value  = "Fruta"
$ws.Cells.Item(177, 7).Value  = 100104
$ws.Cells.Item(177, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(177, 9).Value  = 100104003
$ws.Cells.Item(177, 10).Value = "Membrillo"
$ws.Cells.Item(177, 11).Value = "Champion"
$ws.Cells.Item(177, 12).Value = "Primera"
$ws.Cells.Item(177, 13).Value = 55
$ws.Cells.Item(177, 14).Value = 10000
$ws.Cells.Item(177, 15).Value = 10000
$ws.Cells.Item(177, 16).Value = 10000
$ws.Cells.Item(177, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(177, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(177, 19).Value = 556
$ws.Cells.Item(177, 20).Value = 18
